# The presentation's design/theme ("Integral") color palette is reset back
# to the stock Office "Office Theme" colors -- the same effect as picking
# the default Office color set from the Design > Variants > Colors gallery.
#
# ThemeColorScheme.Item(n).RGB is the read/write surface PowerPoint exposes
# for the 12 DrawingML theme colors; index order follows the standard
# a:clrScheme child order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1..accent6, 11 hlink, 12 folHlink

$p = $ppt.ActivePresentation
$scheme = $p.SlideMaster.Theme.ThemeColorScheme

$officeThemeColors = @(
    0,         # 1  dk1      000000
    16777215,  # 2  lt1      FFFFFF
    6968388,   # 3  dk2      44546A
    15132391,  # 4  lt2      E7E6E6
    13998939,  # 5  accent1  5B9BD5
    3243501,   # 6  accent2  ED7D31
    10855845,  # 7  accent3  A5A5A5
    49407,     # 8  accent4  FFC000
    12874308,  # 9  accent5  4472C4
    4697456,   # 10 accent6  70AD47
    12673797,  # 11 hlink    0563C1
    7491477    # 12 folHlink 954F72
)

for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $scheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
